$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add header labels (previously empty cells)
$ws.Range("A2").Value = "ID"
$ws.Range("B2").Value = "PLATE"

# Row 3: A3 becomes numeric 1 (was text "1"); B3 text changes from "hola" to "dyfu"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "dyfu"

# Row 4: A4 becomes numeric 2 (was text "2"); B4 stays "adios"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "adios"
